$d = $word.ActiveDocument

# The sentence "...payment methods online" was followed by a parenthetical
# hyperlink reference " (https://www.mollie.com/payments)". Find that
# hyperlink and remove it along with its enclosing " (" and ")" so the
# sentence again ends directly in a comma, e.g. "...online, ".
$target = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.Address -eq "https://www.mollie.com/payments") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $linkStart = $target.Range.Start
    $linkEnd = $target.Range.End

    $before = $d.Range($linkStart - 2, $linkStart).Text
    $after = $d.Range($linkEnd, $linkEnd + 1).Text

    $delStart = $linkStart
    $delEnd = $linkEnd
    if ($before -eq " (") { $delStart = $linkStart - 2 }
    if ($after -eq ")") { $delEnd = $linkEnd + 1 }

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
